# Applies the scheduled-runner profit recalculation to the Leve profit
# tables (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N)
# on the affected rows of each job sheet.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 132
$ws.Range("H132").Value = 7411420
$ws.Range("I132").Value = 10003041
$ws.Range("K132").Value = 30009123
$ws.Range("M132").Value = -30006593
# Row 137
$ws.Range("H137").Value = 8116.7075
$ws.Range("J137").Value = 11723.533
$ws.Range("L137").Value = 35170.599
$ws.Range("N137").Value = -40270.599
# Row 138
$ws.Range("H138").Value = 1768.7858
$ws.Range("I138").Value = 1064.75
$ws.Range("K138").Value = 3194.25
$ws.Range("M138").Value = 1945.75

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 28907.555
$ws.Range("J37").Value = 28907.555
$ws.Range("L37").Value = 28907.555
$ws.Range("N37").Value = -29453.555
# Row 44
$ws.Range("H44").Value = 36679
$ws.Range("J44").Value = 36679
$ws.Range("L44").Value = 36679
$ws.Range("N44").Value = -37655
# Row 55
$ws.Range("H55").Value = 44195.25
$ws.Range("J55").Value = 44195.25
$ws.Range("L55").Value = 44195.25
$ws.Range("N55").Value = -44825.25
# Row 61
$ws.Range("H61").Value = 6572.5884
$ws.Range("I61").Value = 3559.3
$ws.Range("K61").Value = 3559.3
$ws.Range("M61").Value = -3347.3
# Row 80
$ws.Range("H80").Value = 46063.8
$ws.Range("J80").Value = 46063.8
$ws.Range("L80").Value = 46063.8
$ws.Range("N80").Value = -48059.8
# Row 83
$ws.Range("H83").Value = 46063.8
$ws.Range("J83").Value = 46063.8
$ws.Range("L83").Value = 138191.4
$ws.Range("N83").Value = -148175.4
# Row 97
$ws.Range("H97").Value = 2850.7368
$ws.Range("I97").Value = 2580.2354
$ws.Range("J97").Value = 5150
$ws.Range("K97").Value = 2580.2354
$ws.Range("L97").Value = 5150
$ws.Range("M97").Value = -2084.2354
$ws.Range("N97").Value = -6142
# Row 132
$ws.Range("H132").Value = 7315.7856
$ws.Range("I132").Value = 5920.7593
$ws.Range("K132").Value = 17762.2779
$ws.Range("M132").Value = -15232.2779
# Row 136
$ws.Range("H136").Value = 6572.5884
$ws.Range("I136").Value = 3559.3
$ws.Range("K136").Value = 10677.9
$ws.Range("M136").Value = -8127.900000000001
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
# Row 139
$ws.Range("H139").Value = 204974.5
$ws.Range("I139").Value = 59950
$ws.Range("J139").Value = 349999
$ws.Range("K139").Value = 59950
$ws.Range("L139").Value = 349999
$ws.Range("M139").Value = -54810
$ws.Range("N139").Value = -360279

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 95
$ws.Range("H95").Value = 22499
$ws.Range("J95").Value = 22499
$ws.Range("L95").Value = 22499
$ws.Range("N95").Value = -27991
# Row 99
$ws.Range("H99").Value = 46681.895
$ws.Range("I99").Value = 119154.336
$ws.Range("K99").Value = 119154.336
$ws.Range("M99").Value = -117656.336
# Row 107
$ws.Range("H107").Value = 521
$ws.Range("I107").Value = 503.0909
$ws.Range("J107").Value = 570.25
$ws.Range("K107").Value = 503.0909
$ws.Range("L107").Value = 570.25
$ws.Range("M107").Value = 1416.9091
$ws.Range("N107").Value = -4410.25
# Row 126
$ws.Range("H126").Value = 46681.895
$ws.Range("I126").Value = 119154.336
$ws.Range("K126").Value = 357463.008
$ws.Range("M126").Value = -354993.008
# Row 132
$ws.Range("H132").Value = 29189.143
$ws.Range("I132").Value = 24037.223
$ws.Range("K132").Value = 72111.66900000001
$ws.Range("M132").Value = -69581.66900000001

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 184
$ws.Range("I2").Value = 165.14285
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 990.8571000000001
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -877.8571000000001
$ws.Range("N2").Value = -1726
# Row 104
$ws.Range("H104").Value = 5870.9287
$ws.Range("J104").Value = 7388.778
$ws.Range("L104").Value = 22166.334
$ws.Range("N104").Value = -27408.334

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 28166.334
$ws.Range("I40").Value = 27249.5
$ws.Range("K40").Value = 27249.5
$ws.Range("M40").Value = -27098.5
# Row 46
$ws.Range("H46").Value = 38126.625
$ws.Range("J46").Value = 52994.4
$ws.Range("L46").Value = 52994.4
$ws.Range("N46").Value = -53306.4
# Row 47
$ws.Range("H47").Value = 32966
$ws.Range("I47").Value = 39449.5
$ws.Range("K47").Value = 39449.5
$ws.Range("M47").Value = -38881.5
# Row 57
$ws.Range("H57").Value = 37515.312
$ws.Range("I57").Value = 40055
$ws.Range("J57").Value = 37346
$ws.Range("K57").Value = 40055
$ws.Range("L57").Value = 37346
$ws.Range("M57").Value = -39235
$ws.Range("N57").Value = -38986
# Row 80
$ws.Range("H80").Value = 9752.462
$ws.Range("I80").Value = 6130.6665
$ws.Range("J80").Value = 12856.857
$ws.Range("K80").Value = 6130.6665
$ws.Range("L80").Value = 12856.857
$ws.Range("M80").Value = -5132.6665
$ws.Range("N80").Value = -14852.857
# Row 83
$ws.Range("H83").Value = 9752.462
$ws.Range("I83").Value = 6130.6665
$ws.Range("J83").Value = 12856.857
$ws.Range("K83").Value = 30653.3325
$ws.Range("L83").Value = 64284.285
$ws.Range("M83").Value = -25661.3325
$ws.Range("N83").Value = -74268.285

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2983.3333
$ws.Range("I7").Value = 2983.3333
$ws.Range("K7").Value = 2983.3333
$ws.Range("M7").Value = -2871.3333
# Row 22
$ws.Range("H22").Value = 1001.12
$ws.Range("I22").Value = 1165.3334
$ws.Range("J22").Value = 754.8
$ws.Range("K22").Value = 1165.3334
$ws.Range("L22").Value = 754.8
$ws.Range("M22").Value = -870.3334
$ws.Range("N22").Value = -1344.8
# Row 27
$ws.Range("H27").Value = 1001.12
$ws.Range("I27").Value = 1165.3334
$ws.Range("J27").Value = 754.8
$ws.Range("K27").Value = 1165.3334
$ws.Range("L27").Value = 754.8
$ws.Range("M27").Value = -1058.3334
$ws.Range("N27").Value = -968.8
# Row 40
$ws.Range("H40").Value = 2950.75
$ws.Range("I40").Value = 2950.75
$ws.Range("K40").Value = 2950.75
$ws.Range("M40").Value = -2814.75
# Row 46
$ws.Range("H46").Value = 1611.92
$ws.Range("I46").Value = 1206.5
$ws.Range("K46").Value = 1206.5
$ws.Range("M46").Value = -1018.5
# Row 122
$ws.Range("H122").Value = 4323.1
$ws.Range("I122").Value = 3489.25
$ws.Range("J122").Value = 4879
$ws.Range("K122").Value = 10467.75
$ws.Range("L122").Value = 14637
$ws.Range("M122").Value = -8017.75
$ws.Range("N122").Value = -19537
# Row 126
$ws.Range("H126").Value = 2983.3333
$ws.Range("I126").Value = 2983.3333
$ws.Range("K126").Value = 8949.999899999999
$ws.Range("M126").Value = -6479.999899999999
# Row 132
$ws.Range("H132").Value = 3591800
$ws.Range("I132").Value = 5296278.5
$ws.Range("K132").Value = 15888835.5
$ws.Range("M132").Value = -15886305.5

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 6818.846
$ws.Range("I126").Value = 7185
$ws.Range("J126").Value = 5598.3335
$ws.Range("K126").Value = 21555
$ws.Range("L126").Value = 16795.0005
$ws.Range("M126").Value = -19085
$ws.Range("N126").Value = -21735.0005
# Row 132
$ws.Range("H132").Value = 13126.574
$ws.Range("I132").Value = 10288.605
$ws.Range("K132").Value = 30865.815
$ws.Range("M132").Value = -28335.815
# Row 136
$ws.Range("H136").Value = 2455.0205
$ws.Range("I136").Value = 2188.3333
$ws.Range("J136").Value = 4055.1428
$ws.Range("K136").Value = 6564.999899999999
$ws.Range("L136").Value = 12165.4284
$ws.Range("M136").Value = -4014.999899999999
$ws.Range("N136").Value = -17265.4284
